$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Room (D) column values
$ws.Range("D2").Value = "R3 (60)"
$ws.Range("D3").Value = "R2 (60)"
$ws.Range("D4").Value = "R7 (60)"
$ws.Range("D5").Value = "R5 (60)"
$ws.Range("D6").Value = "R2 (60)"
$ws.Range("D7").Value = "R6 (60)"
$ws.Range("D8").Value = "R6 (60)"
$ws.Range("D9").Value = "R2 (60)"
$ws.Range("D10").Value = "R3 (60)"
$ws.Range("D13").Value = "R7 (60)"
$ws.Range("D14").Value = "R3 (60)"
$ws.Range("D16").Value = "R5 (60)"
$ws.Range("D17").Value = "R1 (60)"
$ws.Range("D18").Value = "R6 (60)"
$ws.Range("D19").Value = "R7 (60)"

# Update Timing (F) column values
$ws.Range("F2").Value = "Monday 11:00-11:55, Wednesday 10:00-10:55, Thursday 9:00-9:55 (C)"
$ws.Range("F3").Value = "Monday 14:30-15:55, Thursday 16:00-17:25 (P)"
$ws.Range("F4").Value = "Tuesday 14:30-15:55, Friday 16:00-17:25 (R)"
$ws.Range("F5").Value = "Monday 10:00-10:55, Wednesday 9:00-9:55, Friday 9:00-9:55 (B)"
$ws.Range("F6").Value = "Monday 12:00-12:55, Tuesday 9:00-9:55, Friday 11:00-11:55 (D)"
$ws.Range("F7").Value = "Tuesday 10:00-10:55, Thursday 11:00-11:55, Friday 10:00-10:55 (F)"
$ws.Range("F8").Value = "Tuesday 14:30-15:55, Friday 16:00-17:25 (R)"
$ws.Range("F9").Value = "Monday 11:00-11:55, Wednesday 10:00-10:55, Thursday 9:00-9:55 (C)"
$ws.Range("F10").Value = "Tuesday 11:00-11:55,  Wednesday 16:00-17:25, Thursday 12-12:55 (E)"
$ws.Range("F11").Value = "Tuesday 10:00-10:55, Thursday 11:00-11:55, Friday 10:00-10:55 (F)"
$ws.Range("F12").Value = "Monday 12:00-12:55, Tuesday 9:00-9:55, Friday 11:00-11:55 (D)"
$ws.Range("F15").Value = "Tuesday 10:00-10:55, Thursday 11:00-11:55, Friday 10:00-10:55 (F)"
$ws.Range("F17").Value = "Monday 9:00-9:55, Wednesday 11:00-11:55, Thursday 10:00-10:55 (A)"
$ws.Range("F18").Value = "Monday 10:00-10:55, Wednesday 9:00-9:55, Friday 9:00-9:55 (B)"
$ws.Range("F19").Value = "Monday 14:30-15:55, Thursday 16:00-17:25 (P)"
